$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Insert a new column before column A, shifting existing data (A:E -> B:F)
$ws.Columns.Item(1).Insert()

# New header and method names for column A (typed in column order first)
$ws.Range("A1").Value = "Metodo"
$ws.Range("A2").Value = "SMARTER"
$ws.Range("A3").Value = "Fuzzy"
$ws.Range("A4").Value = "TOPSIS"
$ws.Range("A5").Value = "GRA"
$ws.Range("A6").Value = "CODAS"
$ws.Range("A7").Value = "MABAC"
$ws.Range("A8").Value = "VIKOR"
$ws.Range("A9").Value = "PROMETHEE II"

# Updated headers for the shifted columns
$ws.Range("B1").Value = "Rx"
$ws.Range("C1").Value = "Ry"
$ws.Range("D1").Value = "CL"
$ws.Range("E1").Value = "Entropia"
$ws.Range("F1").Value = "SSIM"

# Adjust column widths as Excel would after the content changed
$ws.Columns.Item(1).ColumnWidth = 12.46
$ws.Columns.Item(2).ColumnWidth = 3.2
$ws.Columns.Item(3).ColumnWidth = 2.3
